# Update metrics table with new values from the retrained model ("novo lm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same set applied to every data row, columns B:Q)
$values = @{
    "B" = 0.9999969141016266
    "C" = 0.9990169585194482
    "D" = 0.9999543244690766
    "E" = 0.9999999002982101
    "F" = 0.9999901286903975
    "G" = [double]"2.88055067624043e-06"
    "H" = 0.0009176260715613002
    "I" = [double]"1.315165014746634e-05"
    "J" = [double]"5.003254555116716e-08"
    "K" = [double]"6.600841346508751e-06"
    "L" = [double]"9.999963621106398e-05"
    "M" = 0.00169721851163615
    "N" = 0.9999753128130129
    "O" = 0.001769472617403162
    "P" = 67.51505815013246
    "Q" = 93.11145047236465
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
